# Rename transcript speaker labels in column D ("Speaker") of the Data sheet.
#   "RBD"       -> "T"
#   "Student "  -> "S "   (trailing-space variant, used for unnamed students)
#   "Student 1" -> "S 1"
#   "Student 2" -> "S 2"
# All other speaker names (Michael, Jeff, Brian, Milin, Ankur, Romina, ...)
# are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$lastRow = $ws.Cells.Item($ws.Rows.Count, 4).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    $val = $cell.Value2

    if ($null -eq $val) { continue }

    switch ($val) {
        "RBD"       { $cell.Value2 = "T";   continue }
        "Student "  { $cell.Value2 = "S ";  continue }
        "Student 1" { $cell.Value2 = "S 1"; continue }
        "Student 2" { $cell.Value2 = "S 2"; continue }
    }
}
